$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row => column letter => new value
$updates = @{
    2  = @{ E = 57 }
    3  = @{ E = 27 }
    18 = @{ E = 123 }
    36 = @{ E = 113 }
    43 = @{ E = 28; F = 16; H = 19 }
    48 = @{ E = 35 }
    49 = @{ E = 76; F = 43; H = 60 }
    50 = @{ E = 27 }
    56 = @{ E = 10 }
    59 = @{ E = 11 }
    61 = @{ E = 32; F = 13; H = 23 }
    62 = @{ E = 49; F = 13; H = 27 }
    63 = @{ E = 44; F = 17; H = 25 }
    65 = @{ E = 37 }
    73 = @{ E = 33; F = 12; H = 24 }
    79 = @{ E = 42 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $cols[$col]
    }
}

$wb.Save()
